{"js": "// The worksheet table holds blocks of 5 answer cells every 4th row\n// (rows 0, 4, 8, 12, 16 - zero based) out of 20 rows x 5 columns.\n// Each entry below is [rowIndex, colIndex, expectedOldText, newText].\nconst replacements = [\n  [0, 0, \"48\u00f76=8, 0\", \"38\u00f79=4, 2\"],\n  [0, 1, \"10\u00f79=1, 1\", \"23\u00f78=2, 7\"],\n  [0, 2, \"93\u00f74=23, 1\", \"44\u00f74=11, 0\"],\n  [0, 3, \"90\u00f74=22, 2\", \"65\u00f72=32, 1\"],\n  [0, 4, \"11\u00f72=5, 1\", \"18\u00f75=3, 3\"],\n\n  [4, 0, \"34\u00f76=5, 4\", \"88\u00f77=12, 4\"],\n  [4, 1, \"75\u00f72=37, 1\", \"31\u00f74=7, 3\"],\n  [4, 2, \"35\u00f76=5, 5\", \"46\u00f74=11, 2\"],\n  [4, 3, \"45\u00f76=7, 3\", \"30\u00f75=6, 0\"],\n  [4, 4, \"25\u00f76=4, 1\", \"29\u00f79=3, 2\"],\n\n  [8, 0, \"60\u00f73=20, 0\", \"73\u00f79=8, 1\"],\n  [8, 1, \"51\u00f72=25, 1\", \"77\u00f73=25, 2\"],\n  [8, 2, \"88\u00f75=17, 3\", \"73\u00f77=10, 3\"],\n  [8, 3, \"30\u00f75=6, 0\", \"87\u00f78=10, 7\"],\n  [8, 4, \"36\u00f77=5, 1\", \"82\u00f77=11, 5\"],\n\n  [12, 0, \"70\u00f74=17, 2\", \"78\u00f76=13, 0\"],\n  [12, 1, \"99\u00f77=14, 1\", \"20\u00f78=2, 4\"],\n  [12, 2, \"55\u00f73=18, 1\", \"74\u00f76=12, 2\"],\n  [12, 3, \"26\u00f76=4, 2\", \"90\u00f72=45, 0\"],\n  [12, 4, \"51\u00f79=5, 6\", \"76\u00f72=38, 0\"],\n\n  [16, 0, \"47\u00f73=15, 2\", \"63\u00f74=15, 3\"],\n  [16, 1, \"49\u00f76=8, 1\", \"29\u00f74=7, 1\"],\n  [16, 2, \"60\u00f79=6, 6\", \"97\u00f74=24, 1\"],\n  [16, 3, \"51\u00f76=8, 3\", \"40\u00f76=6, 4\"],\n  [16, 4, \"89\u00f77=12, 5\", \"81\u00f75=16, 1\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Update each target cell in place (getCell uses 0-based row/col indices).\n// Setting `.value` directly rewrites only the cell's text run, leaving the\n// existing run/paragraph formatting (font, size, alignment) untouched.\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  cell.load(\"value\");\n  await context.sync();\n\n  if (cell.value === oldText) {\n    cell.value = newText;\n  } else if (cell.value !== newText) {\n    // Fallback: text didn't match what we expected (unexpected prior edit) -\n    // still try to find/replace it directly within the cell body.\n    const found = cell.body.search(oldText, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    if (found.items.length > 0) {\n      found.items[0].insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# The table holds blocks of 5 answer cells every 4th row (rows 1, 5, 9,\n# 13, 17 - 1 based, as Word COM indexes Cell()) out of 20 rows x 5 columns.\n# Each entry is: row, column, expected old text, new text.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, '48\u00f76=8, 0', '38\u00f79=4, 2'),\n    @(1, 2, '10\u00f79=1, 1', '23\u00f78=2, 7'),\n    @(1, 3, '93\u00f74=23, 1', '44\u00f74=11, 0'),\n    @(1, 4, '90\u00f74=22, 2', '65\u00f72=32, 1'),\n    @(1, 5, '11\u00f72=5, 1', '18\u00f75=3, 3'),\n\n    @(5, 1, '34\u00f76=5, 4', '88\u00f77=12, 4'),\n    @(5, 2, '75\u00f72=37, 1', '31\u00f74=7, 3'),\n    @(5, 3, '35\u00f76=5, 5', '46\u00f74=11, 2'),\n    @(5, 4, '45\u00f76=7, 3', '30\u00f75=6, 0'),\n    @(5, 5, '25\u00f76=4, 1', '29\u00f79=3, 2'),\n\n    @(9, 1, '60\u00f73=20, 0', '73\u00f79=8, 1'),\n    @(9, 2, '51\u00f72=25, 1', '77\u00f73=25, 2'),\n    @(9, 3, '88\u00f75=17, 3', '73\u00f77=10, 3'),\n    @(9, 4, '30\u00f75=6, 0', '87\u00f78=10, 7'),\n    @(9, 5, '36\u00f77=5, 1', '82\u00f77=11, 5'),\n\n    @(13, 1, '70\u00f74=17, 2', '78\u00f76=13, 0'),\n    @(13, 2, '99\u00f77=14, 1', '20\u00f78=2, 4'),\n    @(13, 3, '55\u00f73=18, 1', '74\u00f76=12, 2'),\n    @(13, 4, '26\u00f76=4, 2', '90\u00f72=45, 0'),\n    @(13, 5, '51\u00f79=5, 6', '76\u00f72=38, 0'),\n\n    @(17, 1, '47\u00f73=15, 2', '63\u00f74=15, 3'),\n    @(17, 2, '49\u00f76=8, 1', '29\u00f74=7, 1'),\n    @(17, 3, '60\u00f79=6, 6', '97\u00f74=24, 1'),\n    @(17, 4, '51\u00f76=8, 3', '40\u00f76=6, 4'),\n    @(17, 5, '89\u00f77=12, 5', '81\u00f75=16, 1')\n)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $oldText = $entry[2]\n    $newText = $entry[3]\n\n    $cell = $t.Cell($row, $col)\n    # Cell.Range.Text includes the trailing cell-mark (\\r\\a); strip it off\n    # before comparing against the expected original text.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -eq $oldText -or $current -ne $newText) {\n        # Assigning Range.Text replaces just the cell's text run while\n        # keeping the run/paragraph formatting (font, size, alignment).\n        $cell.Range.Text = $newText\n    }\n}\n"}
